$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "dry hands tennis"
$ws.Range("A2").Value = "dry inhaler"
$ws.Range("A3").Value = "dry skin leggings"
$ws.Range("A4").Value = "duty gear stand"
$ws.Range("A5").Value = "dynamic athletica womens"
$ws.Range("A6").Value = "dynamic gear"
$ws.Range("A7").Value = "dynamic muscle"
$ws.Range("A8").Value = "dynamic ring"
$ws.Range("A9").Value = "dynamic tape 2 inch"
$ws.Range("A10").Value = "effect apparel"
$ws.Range("A11").Value = "elastic joint"
$ws.Range("A12").Value = "elastic tights for women"
$ws.Range("A13").Value = "element fitness"
$ws.Range("A14").Value = "element pants"
$ws.Range("A15").Value = "elements cigarette"
$ws.Range("A16").Value = "ella miles"
$ws.Range("A17").Value = "endo apparel"
$ws.Range("A18").Value = "endurance 2"
$ws.Range("A19").Value = "endurance 360"
$ws.Range("A20").Value = "endurance a year"
$ws.Range("A21").Value = "endurance bands"
$ws.Range("A22").Value = "endurance bike men"
$ws.Range("A23").Value = "endurance blend"
$ws.Range("A24").Value = "endurance equipment"
$ws.Range("A25").Value = "endurance generator"
$ws.Range("A26").Value = "endurance generator tights"
$ws.Range("A27").Value = "endurance plus"
$ws.Range("A28").Value = "endurance recovery"
$ws.Range("A29").Value = "endurance riding"
$ws.Range("A30").Value = "endurance run"
$ws.Range("A31").Value = "endurance running"
$ws.Range("A32").Value = "endurance training equipment"
$ws.Range("A33").Value = "equestrian clothes women"
$ws.Range("A34").Value = "equestrian winter gear"
$ws.Range("A35").Value = "equipment for marijuana"
$ws.Range("A36").Value = "es 130"
$ws.Range("A37").Value = "essential brace knee brace"
$ws.Range("A38").Value = "essential knee brace"
$ws.Range("A39").Value = "essential knee brace spring"
$ws.Range("A40").Value = "essential pocket chart"
$ws.Range("A41").Value = "essentials performance logo pants"
$ws.Range("A42").Value = "eve cigarettes"
$ws.Range("A43").Value = "even faster websites"
$ws.Range("A44").Value = "ever after winter"
$ws.Range("A45").Value = "ever brand leggings"
$ws.Range("A46").Value = "ever wear"
$ws.Range("A47").Value = "everyday knee brace"
$ws.Range("A48").Value = "everyday leggings for women"
$ws.Range("A49").Value = "everyday skirts"
$ws.Range("A50").Value = "everyday wear leggings"
$ws.Range("A51").Value = "excersise clothes"
$ws.Range("A52").Value = "excersize capris"
$ws.Range("A53").Value = "exchange items purchased"
$ws.Range("A54").Value = "exercise & fitness apparel"
$ws.Range("A55").Value = "exercise ball flat bottom"
$ws.Range("A56").Value = "exercise ball with base and bands"
$ws.Range("A57").Value = "exercise bike generator"
$ws.Range("A58").Value = "exercise capris"
$ws.Range("A59").Value = "exercise charts for women"
$ws.Range("A60").Value = "exercise clothes"
$ws.Range("A61").Value = "exercise clothing for women plus size"
$ws.Range("A62").Value = "exercise compression leggings"
$ws.Range("A63").Value = "exercise equipment lower back"
$ws.Range("A64").Value = "exercise for stomach"
$ws.Range("A65").Value = "exercise for thighs"
$ws.Range("A66").Value = "exercise glutes"
$ws.Range("A67").Value = "exercise jumpsuit"
$ws.Range("A68").Value = "exercise knee support"
$ws.Range("A69").Value = "exercise leggings"
$ws.Range("A70").Value = "exercise leggings for women"
$ws.Range("A71").Value = "exercise leggings for women plus size"
$ws.Range("A72").Value = "exercise leggings plus size"
$ws.Range("A73").Value = "exercise pants compression women"
$ws.Range("A74").Value = "exercise pants for plus size women"
$ws.Range("A75").Value = "exercise pants for women high waist"
$ws.Range("A76").Value = "exercise plus size"
$ws.Range("A77").Value = "exercise ring squeeze"
$ws.Range("A78").Value = "exercise shorts plus size"
$ws.Range("A79").Value = "exercise support"
$ws.Range("A80").Value = "exercise thigh"
$ws.Range("A81").Value = "exercise thighs"
$ws.Range("A82").Value = "exercise tights"
$ws.Range("A83").Value = "exercise tights reflective women pocket"
$ws.Range("A84").Value = "exercise waist"
$ws.Range("A85").Value = "exercise wear for women"
$ws.Range("A86").Value = "exercises clothes for women"
$ws.Range("A87").Value = "exo 13"
$ws.Range("A88").Value = "exo apparel"
$ws.Range("A89").Value = "exo clothes"
$ws.Range("A90").Value = "exo exact"
$ws.Range("A91").Value = "exo force"
$ws.Range("A92").Value = "exo gear"
$ws.Range("A93").Value = "exo international"
$ws.Range("A94").Value = "exo jackets"
$ws.Range("A95").Value = "exo knee brace"
$ws.Range("A96").Value = "exo knee sleeves"
$ws.Range("A97").Value = "exo leggings"
$ws.Range("A98").Value = "exo note"
$ws.Range("A99").Value = "exo pants"
$ws.Range("A100").Value = "exo patch"
